$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows (107-110), columns A-P
$data = @(
    @(44351, 600, 3152, 911, 6879, 524, 2195, 0, 28500, 0, 0, 0, 3674, 13, 1087, 44400),
    @(44354, 902, 3423, 572, 6258, 1066, 2293, 1000, 28500, 0, 0, 0, 3674, 278, 1008, 44148),
    @(44355, 400, 3523, 854, 6275, 212, 2291, 0, 28500, 0, 0, 0, 3674, 122, 1008, 44263),
    @(44356, 930, 3698, 700, 6275, 75, 2316, 0, 28500, 0, 0, 0, 3674, 77, 999, 44463)
)

$startRow = 107
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowValues = $data[$i]
    for ($c = 1; $c -le $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
    # Column A gets the date number format matching the rest of the date column
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd"
}
